$d = $word.ActiveDocument

# Find the two "target" paragraphs that carry the "De asemenea, sunt de acord ca
# Facultatea ..." consent statement, and the paragraph just before it, plus the
# "Prin prezenta declar ..." paragraph right above - all identified by their text.
$consentPara = $null
$blankBeforePara = $null
$introPara = $null

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $t = $para.Range.Text
    if ($t.StartsWith("De asemenea, sunt de acord ca Facultatea")) {
        $consentPara = $para
        $blankBeforePara = $d.Paragraphs.Item($i - 1)
        $introPara = $d.Paragraphs.Item($i - 2)
    }
}

# The home-page "cover" bookmark Word drops at the last edit location - it used to
# sit right after the page break leading into the consent-declaration page.
# Move it onto the new typed text further down (see step below); first drop the
# old one so the `page15` bookmark can re-claim id 2.
$d.Bookmarks("_GoBack").Delete()

# New first-line indents on the intro paragraph, the blank spacer paragraph, and
# the consent paragraph itself, matching the rest of the document's body style.
$introPara.Range.ParagraphFormat.FirstLineIndent = 36
$blankBeforePara.Range.ParagraphFormat.FirstLineIndent = 36
$consentPara.Range.ParagraphFormat.FirstLineIndent = 36

# Re-drop "_GoBack" at the point the author last left the cursor: right after
# "De a" within the consent paragraph, splitting the run in two.
$insertPos = $consentPara.Range.Start + 4
$r = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $r)
